$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-16 20:48:54'
$ws.Range("I2").Value = '21.2 mm'
$ws.Range("E3").Value = '2026-02-16 20:48:57'
$ws.Range("G3").Value = '243 cm'
$ws.Range("I3").Value = '10.1 mm'
$ws.Range("N3").Value = '-2.4 °C 20:24 TU'
$ws.Range("E4").Value = '2026-02-16 20:49:00'
$ws.Range("O4").Value = '13.6 °C'
$ws.Range("E5").Value = '2026-02-16 20:49:02'
$ws.Range("I5").Value = '24.9 mm'
$ws.Range("L5").Value = '50.8 km/h - 336º 20:21 TU'
$ws.Range("N5").Value = '-2.0 °C 20:21 TU'
$ws.Range("E6").Value = '2026-02-16 20:49:05'
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = '67%'
$ws.Range("E7").Value = '2026-02-16 20:49:08'
$ws.Range("O7").Value = '16.5 °C'
$ws.Range("E8").Value = '2026-02-16 20:49:11'
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = '66%'
$ws.Range("J8").Value = '1012.9 hPa'
$ws.Range("E9").Value = '2026-02-16 20:49:13'
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = '70%'
$ws.Range("E10").Value = '2026-02-16 20:49:16'
$ws.Range("O10").Value = '10.9 °C'
$ws.Range("E11").Value = '2026-02-16 20:49:19'
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = '77%'
$ws.Range("E12").Value = '2026-02-16 20:49:21'
$ws.Range("O12").Value = '11.0 °C'
$ws.Range("E13").Value = '2026-02-16 20:49:24'
$ws.Range("E14").Value = '2026-02-16 20:49:27'
$ws.Range("E15").Value = '2026-02-16 20:49:28'
$ws.Range("E16").Value = '2026-02-16 20:49:29'
$ws.Range("N16").Value = '-2.4 °C 20:06 TU'
$ws.Range("O16").Value = '-0.2 °C'
$ws.Range("E17").Value = '2026-02-16 20:49:30'
$ws.Range("O17").Value = '6.0 °C'
$ws.Range("E18").Value = '2026-02-16 20:49:31'
$ws.Range("J18").Value = '1012.7 hPa'
$ws.Range("O18").Value = '10.8 °C'
$ws.Range("E19").Value = '2026-02-16 20:49:33'
$ws.Range("E20").Value = '2026-02-16 20:49:34'
$ws.Range("I20").Value = '0.6 mm'
$ws.Range("E21").Value = '2026-02-16 20:49:35'
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = '69%'
$ws.Range("J21").Value = '1014.1 hPa'
$ws.Range("L21").Value = '50.8 km/h - 285º 20:27 TU'
$ws.Range("O21").Value = '8.7 °C'
$ws.Range("E22").Value = '2026-02-16 20:49:36'
$ws.Range("E23").Value = '2026-02-16 20:49:37'
$ws.Range("I23").Value = '15.1 mm'
$ws.Range("N23").Value = '-2.8 °C 20:29 TU'
$ws.Range("O23").Value = '-0.7 °C'
$ws.Range("E24").Value = '2026-02-16 20:49:40'
$ws.Range("E25").Value = '2026-02-16 20:49:43'
$ws.Range("I25").Value = '6.3 mm'
$ws.Range("N25").Value = '-0.9 °C 20:10 TU'
$ws.Range("E26").Value = '2026-02-16 20:49:46'
$ws.Range("E27").Value = '2026-02-16 20:49:48'
$ws.Range("E28").Value = '2026-02-16 20:49:51'
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = '72%'
$ws.Range("E29").Value = '2026-02-16 20:49:53'
$ws.Range("E30").Value = '2026-02-16 20:49:56'
$ws.Range("E31").Value = '2026-02-16 20:49:59'
$ws.Range("N31").Value = '12.3 °C 20:29 TU'
$ws.Range("O31").Value = '14.5 °C'
$ws.Range("E32").Value = '2026-02-16 20:50:02'
$ws.Range("O32").Value = '8.7 °C'
$ws.Range("E33").Value = '2026-02-16 20:50:04'
$ws.Range("E34").Value = '2026-02-16 20:50:07'
$ws.Range("E35").Value = '2026-02-16 20:50:10'
$ws.Range("I35").Value = '1.1 mm'
$ws.Range("J35").Value = '1016.6 hPa'
$ws.Range("E36").Value = '2026-02-16 20:50:12'
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = '72%'
$ws.Range("J36").Value = '1012.6 hPa'
$ws.Range("O36").Value = '12.0 °C'
$ws.Range("E37").Value = '2026-02-16 20:50:15'
$ws.Range("E38").Value = '2026-02-16 20:50:18'
$ws.Range("E39").Value = '2026-02-16 20:50:21'
$ws.Range("I39").Value = '4.1 mm'
$ws.Range("N39").Value = '-2.1 °C 20:05 TU'
$ws.Range("E40").Value = '2026-02-16 20:50:23'
$ws.Range("O40").Value = '7.1 °C'
$ws.Range("E41").Value = '2026-02-16 20:50:26'
$ws.Range("E42").Value = '2026-02-16 20:50:28'
$ws.Range("E43").Value = '2026-02-16 20:50:31'
$ws.Range("E44").Value = '2026-02-16 20:50:34'
$ws.Range("I44").Value = '11.3 mm'
$ws.Range("N44").Value = '-2.2 °C 20:29 TU'
$ws.Range("O44").Value = '-0.2 °C'
$ws.Range("E45").Value = '2026-02-16 20:50:37'
$ws.Range("I45").Value = '17.9 mm'
$ws.Range("E46").Value = '2026-02-16 20:50:39'
